$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 721

# Update row 3 values
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 281

# Remove row 4 entirely (shifts nothing up since it's the last row)
$ws.Rows.Item(4).Delete()
